$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Step 1: Insert a new row at 18. This shifts the old row 19 header
#          "Two pointers (Basics)" down to row 20, and the old row 21
#          "Valid Palindrome" data row down to row 22, preserving the
#          blank spacer rows in between (19 and 21).
# ------------------------------------------------------------------
$ws.Rows("18:18").Insert()

# New row 18 data ("Tuple with Same Product") - continuation of the
# "Arrays(Advance)" section, mirrors the style of rows 16/17 (B col
# uses style 4, date column uses the same date number format).
$ws.Range("A18").Value = 1726
$ws.Range("B18").Value = "Tuple with Same Product"
$ws.Range("C18").Value = "Medium"
$ws.Range("D18").Value = "Arrays,hashmap,Counting,Combinations"
$ws.Range("E18").Value = 45694
$ws.Range("E18").NumberFormat = $ws.Range("E17").NumberFormat

# ------------------------------------------------------------------
# Step 2: Insert a new row at 23 (a plain blank insert below current
#          row 22) to hold the new "Two Sum II" entry.
# ------------------------------------------------------------------
$ws.Rows("23:23").Insert()

$ws.Range("A23").Value = 167
$ws.Range("B23").Value = "Two Sum II - Input Array Is Sorted"
$ws.Range("C23").Value = "Medium"
$ws.Range("D23").Value = "Arrays,Two pointers,Binary Search"
$ws.Range("E23").Value = 45694
$ws.Range("E23").NumberFormat = $ws.Range("E22").NumberFormat

# ------------------------------------------------------------------
# Step 3: Re-home the "Two pointers (Basics)" section header (row 20)
#          onto the same style used by the other section headers
#          ("Arrays (Basics)" / "Arrays(Advance)").
# ------------------------------------------------------------------
$ws.Range("B3").Copy() | Out-Null
$pasted = $ws.Range("B20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Step 4: Column D width & selection / active cell bookkeeping.
# ------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 41.55

$ws.Range("D24").Select() | Out-Null
